$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 296
$ws1.Range("F3").Value = 300
$ws1.Range("F7").Value = 9781
$ws1.Range("F10").Value = 135
$ws1.Range("F11").Value = 125
$ws1.Range("F17").Value = 281
$ws1.Range("F18").Value = 792
$ws1.Range("F19").Value = 47

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 296
$ws4.Range("F3").Value = 300
$ws4.Range("F7").Value = 9782
$ws4.Range("F10").Value = 135
$ws4.Range("F11").Value = 125
$ws4.Range("F17").Value = 281
$ws4.Range("F18").Value = 792
$ws4.Range("F19").Value = 47
